# Update gh-pages to output generated at 456a3b4
# Apply the refreshed 想去人数 (F) / 最低票价 (G) figures to the "展览" and
# "全部类型" sheets (the two sheets that list exhibition rows).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" -----------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Cells.Item(2, 6).Value  = 15206   # F2
$ws1.Cells.Item(3, 6).Value  = 19620   # F3
$ws1.Cells.Item(5, 6).Value  = 183     # F5
$ws1.Cells.Item(5, 7).Value  = 65      # G5
$ws1.Cells.Item(13, 7).Value = 178     # G13
$ws1.Cells.Item(15, 7).Value = 68      # G15
$ws1.Cells.Item(17, 6).Value = 1536    # F17
$ws1.Cells.Item(20, 6).Value = 121     # F20
$ws1.Cells.Item(22, 6).Value = 8286    # F22
$ws1.Cells.Item(23, 6).Value = 994     # F23
$ws1.Cells.Item(25, 6).Value = 13      # F25
$ws1.Cells.Item(27, 6).Value = 1281    # F27
$ws1.Cells.Item(28, 6).Value = 36      # F28
$ws1.Cells.Item(31, 6).Value = 6617    # F31
$ws1.Cells.Item(32, 6).Value = 140     # F32
$ws1.Cells.Item(33, 6).Value = 81      # F33
$ws1.Cells.Item(36, 6).Value = 319     # F36
$ws1.Cells.Item(37, 6).Value = 5666    # F37
$ws1.Cells.Item(38, 6).Value = 1024    # F38

# --- Sheet "全部类型" ---------------------------------------------------
$ws2 = $wb.Worksheets.Item("全部类型")

$ws2.Cells.Item(2, 6).Value  = 15206   # F2
$ws2.Cells.Item(3, 6).Value  = 19620   # F3
$ws2.Cells.Item(5, 6).Value  = 183     # F5
$ws2.Cells.Item(5, 7).Value  = 65      # G5
$ws2.Cells.Item(13, 7).Value = 178     # G13
$ws2.Cells.Item(15, 7).Value = 68      # G15
$ws2.Cells.Item(17, 6).Value = 1536    # F17
$ws2.Cells.Item(21, 6).Value = 121     # F21
$ws2.Cells.Item(23, 6).Value = 8286    # F23
$ws2.Cells.Item(26, 6).Value = 13      # F26
$ws2.Cells.Item(28, 6).Value = 1281    # F28
$ws2.Cells.Item(29, 6).Value = 36      # F29
$ws2.Cells.Item(34, 6).Value = 6617    # F34
$ws2.Cells.Item(35, 6).Value = 140     # F35
$ws2.Cells.Item(36, 6).Value = 81      # F36
$ws2.Cells.Item(39, 6).Value = 319     # F39
$ws2.Cells.Item(40, 6).Value = 5666    # F40
$ws2.Cells.Item(41, 6).Value = 1024    # F41
